# Actualización 11 de Mayo - Mañana
$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 1P" ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D3").Value = 2
$ws1.Range("E3").Value = 5
$ws1.Range("H3").Value = 6.2

$ws1.Range("D4").Value = 0
$ws1.Range("E4").Value = 6
$ws1.Range("H4").Value = 6.2

# --- Sheet "Estadisticos 2P" ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("D2").Value = 10
$ws2.Range("E2").Value = 9
$ws2.Range("F2").Value = 18
$ws2.Range("G2").Value = 64.29000000000001

$ws2.Range("D3").Value = 11
$ws2.Range("F3").Value = 17
$ws2.Range("G3").Value = 60.71

$ws2.Range("D4").Value = 3
$ws2.Range("E4").Value = 3
$ws2.Range("F4").Value = 19
$ws2.Range("G4").Value = 86.36

# --- Sheet "Estadisticos Final" ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D3").Value = 2
$ws3.Range("E3").Value = 5
$ws3.Range("H3").Value = 6.4

$ws3.Range("D4").Value = 0
$ws3.Range("E4").Value = 3
$ws3.Range("H4").Value = 6.5
